# Pivots.xlsx update — roll the "Nifty" sheet's two live pivot columns (G, H)
# forward by one trading day: the former "H" (latest) day's figures slide
# into "G" (previous day), and brand-new figures are entered into "H".
# All pivot/support-resistance formulas in columns G:H recalculate
# automatically from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nifty")

# ---------------------------------------------------------------------
# 1) Header date row + High/Low/Close inputs (rows 1-4): shift G<-H, new H
# ---------------------------------------------------------------------
$oldH1 = $ws.Range("H1").Value2
$oldH2 = $ws.Range("H2").Value2
$oldH3 = $ws.Range("H3").Value2
$oldH4 = $ws.Range("H4").Value2

$ws.Range("G1").Value = $oldH1
$ws.Range("G2").Value = $oldH2
$ws.Range("G3").Value = $oldH3
$ws.Range("G4").Value = $oldH4

$ws.Range("H1").Value = 43501
$ws.Range("H2").Value = 10956.7
$ws.Range("H3").Value = 10886.7
$ws.Range("H4").Value = 10934.35

# ---------------------------------------------------------------------
# 2) Manually-pasted pivot snapshot rows (39-42, 44-47): the "G" (previous
#    day) cell picks up what used to be in "H" (latest day); "H" itself is
#    left exactly as it was (it will get a fresh paste on the next update).
# ---------------------------------------------------------------------
$snapshotRows = 39, 40, 41, 42, 44, 45, 46, 47
foreach ($r in $snapshotRows) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $gCell = $ws.Cells.Item($r, 7)   # column G
    if ($hCell.Value2 -eq $null) {
        $gCell.ClearContents()
    } else {
        $gCell.Value = $hCell.Value2
    }
}

# ---------------------------------------------------------------------
# 3) Carry the bold "current pivot" highlight on row 33 from H over to G
#    as well (H keeps it too).
# ---------------------------------------------------------------------
$ws.Range("G33").Font.Bold = $true

# ---------------------------------------------------------------------
# 4) Sheet view: drop the frozen "topLeftCell" scroll anchor and move the
#    selection from the old H24 cell onto the newly-entered G1:H4 block.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G1:H4").Select()

# ---------------------------------------------------------------------
# 5) Cosmetic workbook-level metadata (best effort — author's machine path
#    and the saved window size in the OOXML bookViews element).
# ---------------------------------------------------------------------
try { $excel.ActiveWindow.Width = 9504 } catch {}
try { $excel.ActiveWindow.Height = 7488 } catch {}
try { $wb.Path = "F:\F-Anbu\Business\" } catch {}
